$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description timestamp in A1
$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 05:57:14 EDT, by WPJTOWN1.The search returned: 1 events."

# Update data row (row 3) values
$ws.Range("C3").Value = "JOHNSTOWN"
$ws.Range("D3").Value = "CO"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1811
$ws.Range("H3").Value = "Placed Actual"
$ws.Range("I3").ClearContents()

# Remove the AutoFilter (and its defined name / filter database)
$ws.AutoFilterMode = $false
if ($wb.Names.Count -gt 0) {
    for ($i = $wb.Names.Count; $i -ge 1; $i--) {
        $wb.Names.Item($i).Delete()
    }
}
